# Apply updated cryptocurrency price/volume values to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.359.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.595.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.05'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.65%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0601'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.824.33'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.603.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.390.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.83%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.399.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.11%  '
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -5.42%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.60%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.814'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  +8.86%  '
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.984'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.735.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0527'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.35%  '
